$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.5058561382584799
$ws.Range("C2").Value = 0.143701208575493
$ws.Range("E2").Value = 0.2336753527718756
$ws.Range("F2").Value = 2.895405934404607
$ws.Range("G2").Value = 0.002523578501680898
$ws.Range("I2").Value = 1.18932665027797
$ws.Range("J2").Value = 0.1103267142276909
$ws.Range("K2").Value = 0.6388100547983697
$ws.Range("M2").Value = 0.4173734276955585
$ws.Range("N2").Value = 2.500046538005279
# Row 3
$ws.Range("B3").Value = 0.4726270907902119
$ws.Range("C3").Value = 0.1358393148541381
$ws.Range("E3").Value = 0.2289736659987511
$ws.Range("F3").Value = 2.873027262623481
$ws.Range("G3").Value = 0.002527187915113905
$ws.Range("I3").Value = 1.19083112469368
$ws.Range("J3").Value = 0.1110801525659504
$ws.Range("K3").Value = 0.5987489056767288
$ws.Range("M3").Value = 0.4011072879575153
$ws.Range("N3").Value = 2.518167781686913
# Row 4
$ws.Range("B4").Value = 0.4524925900373375
$ws.Range("C4").Value = 0.1311008126469062
$ws.Range("E4").Value = 0.2261997684549897
$ws.Range("F4").Value = 2.860740913906128
$ws.Range("G4").Value = 0.002529520943273767
$ws.Range("I4").Value = 1.192282533542048
$ws.Range("J4").Value = 0.111574886244199
$ws.Range("K4").Value = 0.5745053756218113
$ws.Range("M4").Value = 0.3913505902154881
$ws.Range("N4").Value = 2.530020781158136
# Row 5
$ws.Range("B5").Value = 0.4443551177393203
$ws.Range("C5").Value = 0.1291920457450857
$ws.Range("E5").Value = 0.2250978316777363
$ws.Range("F5").Value = 2.856099602391922
$ws.Range("G5").Value = 0.00253050114300136
$ws.Range("I5").Value = 1.193006556400491
$ws.Range("J5").Value = 0.1117845683227205
$ws.Range("K5").Value = 0.5647149183314752
$ws.Range("M5").Value = 0.3874327190115352
$ws.Range("N5").Value = 2.535033506648929
# Row 6
$ws.Range("B6").Value = 0.443007974861132
$ws.Range("C6").Value = 0.1288764350276494
$ws.Range("E6").Value = 0.2249165755059366
$ws.Range("F6").Value = 2.855350984106295
$ws.Range("G6").Value = 0.00253066568721953
$ws.Range("I6").Value = 1.193134782760268
$ws.Range("J6").Value = 0.1118198734533671
$ws.Range("K6").Value = 0.5630945921380714
$ws.Range("M6").Value = 0.3867856672922017
$ws.Range("N6").Value = 2.535876886106294
# Row 7
$ws.Range("B7").Value = 0.4523825718771661
$ws.Range("C7").Value = 0.1310749805489593
$ws.Range("E7").Value = 0.2261847921007245
$ws.Range("F7").Value = 2.860676840040469
$ws.Range("G7").Value = 0.00252953404312334
$ws.Range("I7").Value = 1.192291761394884
$ws.Range("J7").Value = 0.1115776814017035
$ws.Range("K7").Value = 0.5743729779293858
$ws.Range("M7").Value = 0.3912975172988737
$ws.Range("N7").Value = 2.530087645669482
# Row 8
$ws.Range("B8").Value = 0.4943430962980813
$ws.Range("C8").Value = 0.140971947324303
$ws.Range("E8").Value = 0.2320307917482722
$ws.Range("F8").Value = 2.887387799317139
$ws.Range("G8").Value = 0.002524798834723668
$ws.Range("I8").Value = 1.189735802427371
$ws.Range("J8").Value = 0.1105798317982103
$ws.Range("K8").Value = 0.6249233962926439
$ws.Range("M8").Value = 0.4117169769957556
$ws.Range("N8").Value = 2.506143873016676
# Row 9
$ws.Range("B9").Value = 0.578759275947192
$ws.Range("C9").Value = 0.1610896665013399
$ws.Range("E9").Value = 0.2443899461519479
$ws.Range("F9").Value = 2.951322595395297
$ws.Range("G9").Value = 0.002516435879004613
$ws.Range("I9").Value = 1.188916893592044
$ws.Range("J9").Value = 0.1088780167127865
$ws.Range("K9").Value = 0.7268726128205856
$ws.Range("M9").Value = 0.4535918977838449
$ws.Range("N9").Value = 2.464960826764752
# Row 10
$ws.Range("B10").Value = 0.6420909602639995
$ws.Range("C10").Value = 0.1763125216735659
$ws.Range("E10").Value = 0.2540157637220943
$ws.Range("F10").Value = 3.005371039136747
$ws.Range("G10").Value = 0.002510848199059515
$ws.Range("I10").Value = 1.190882398633647
$ws.Range("J10").Value = 0.1077832711081079
$ws.Range("K10").Value = 0.8035160552725245
$ws.Range("M10").Value = 0.485480227536442
$ws.Range("N10").Value = 2.438229772628226
# Row 11
$ws.Range("B11").Value = 0.671189850534688
$ws.Range("C11").Value = 0.1833360445405106
$ws.Range("E11").Value = 0.2585133143525979
$ws.Range("F11").Value = 3.031502957390757
$ws.Range("G11").Value = 0.002508425804472009
$ws.Range("I11").Value = 1.192336467083869
$ws.Range("J11").Value = 0.1073190603044836
$ws.Range("K11").Value = 0.838766619895523
$ws.Range("M11").Value = 0.5002323650431819
$ws.Range("N11").Value = 2.426836621969628
# Row 12
$ws.Range("B12").Value = 0.6822504879393421
$ws.Range("C12").Value = 0.1860099871155114
$ws.Range("E12").Value = 0.260233468956045
$ws.Range("F12").Value = 3.041621039258843
$ws.Range("G12").Value = 0.002507525589985213
$ws.Range("I12").Value = 1.192967783649458
$ws.Range("J12").Value = 0.1071481396039147
$ws.Range("K12").Value = 0.8521707169772981
$ws.Range("M12").Value = 0.5058540275201651
$ws.Range("N12").Value = 2.422632803564724
# Row 13
$ws.Range("B13").Value = 0.6798665357569575
$ws.Range("C13").Value = 0.1854334685187382
$ws.Range("E13").Value = 0.2598622461383258
$ws.Range("F13").Value = 3.039432028731909
$ws.Range("G13").Value = 0.002507718708494719
$ws.Range("I13").Value = 1.192828226745
$ws.Range("J13").Value = 0.1071847339227947
$ws.Range("K13").Value = 0.8492814360241425
$ws.Range("M13").Value = 0.5046417301489896
$ws.Range("N13").Value = 2.42353325184493
# Row 14
$ws.Range("B14").Value = 0.6720989839261051
$ws.Range("C14").Value = 0.1835557444244102
$ws.Range("E14").Value = 0.2586544914415256
$ws.Range("F14").Value = 3.032330917150205
$ws.Range("G14").Value = 0.002508351401265796
$ws.Range("I14").Value = 1.192386787686942
$ws.Range("J14").Value = 0.1073049010000204
$ws.Range("K14").Value = 0.8398682706776412
$ws.Range("M14").Value = 0.5006941542938108
$ws.Range("N14").Value = 2.426488554940093
# Row 15
$ws.Range("B15").Value = 0.6673465401186718
$ws.Range("C15").Value = 0.1824074481401965
$ws.Range("E15").Value = 0.2579169233699403
$ws.Range("F15").Value = 3.02801026587963
$ws.Range("G15").Value = 0.002508741167606731
$ws.Range("I15").Value = 1.192126907202407
$ws.Range("J15").Value = 0.1073791406716396
$ws.Range("K15").Value = 0.8341096623851172
$ws.Range("M15").Value = 0.4982807528220334
$ws.Range("N15").Value = 2.428313163512001
# Row 16
$ws.Range("B16").Value = 0.6401950913150642
$ws.Range("C16").Value = 0.1758555110172324
$ws.Range("E16").Value = 0.2537242235415036
$ws.Range("F16").Value = 3.003694373655492
$ws.Range("G16").Value = 0.002511008904902912
$ws.Range("I16").Value = 1.190798654139442
$ws.Range("J16").Value = 0.1078142889367317
$ws.Range("K16").Value = 0.8012201051492696
$ws.Range("M16").Value = 0.4845210883304247
$ws.Range("N16").Value = 2.438989794633791
# Row 17
$ws.Range("B17").Value = 0.6236125466846829
$ws.Range("C17").Value = 0.1718614356771582
$ws.Range("E17").Value = 0.2511825146559445
$ws.Range("F17").Value = 2.989173301735917
$ws.Range("G17").Value = 0.002512430625944417
$ws.Range("I17").Value = 1.190127352404097
$ws.Range("J17").Value = 0.1080898982599994
$ws.Range("K17").Value = 0.7811421244305166
$ws.Range("M17").Value = 0.4761429634308101
$ws.Range("N17").Value = 2.44573615984816
# Row 18
$ws.Range("B18").Value = 0.6141019086948347
$ws.Range("C18").Value = 0.1695734202177164
$ws.Range("E18").Value = 0.2497317676854749
$ws.Range("F18").Value = 2.980966569245425
$ws.Range("G18").Value = 0.00251325961228534
$ws.Range("I18").Value = 1.189793932108913
$ws.Range("J18").Value = 0.1082516019597808
$ws.Range("K18").Value = 0.7696300401599672
$ws.Range("M18").Value = 0.4713472413089974
$ws.Range("N18").Value = 2.449688664295486
# Row 19
$ws.Range("B19").Value = 0.6108864434238228
$ws.Range("C19").Value = 0.1688003263919029
$ws.Range("E19").Value = 0.2492424905864823
$ws.Range("F19").Value = 2.978212874260194
$ws.Range("G19").Value = 0.002513542227890755
$ws.Range("I19").Value = 1.189690086892597
$ws.Range("J19").Value = 0.108306898122672
$ws.Range("K19").Value = 0.7657384673790659
$ws.Range("M19").Value = 0.469727469601672
$ws.Range("N19").Value = 2.451039303190186
# Row 20
$ws.Range("B20").Value = 0.6253749717302242
$ws.Range("C20").Value = 0.1722856514867033
$ws.Range("E20").Value = 0.2514519274582909
$ws.Range("F20").Value = 2.990704043549215
$ws.Range("G20").Value = 0.002512278117562201
$ws.Range("I20").Value = 1.19019335866939
$ws.Range("J20").Value = 0.1080602299613354
$ws.Range("K20").Value = 0.7832757084256627
$ws.Range("M20").Value = 0.4770324329774951
$ws.Range("N20").Value = 2.445010526351858
# Row 21
$ws.Range("B21").Value = 0.6743793756496359
$ws.Range("C21").Value = 0.1841068888615496
$ws.Range("E21").Value = 0.259008776418419
$ws.Range("F21").Value = 3.034410644804808
$ws.Range("G21").Value = 0.00250816510094823
$ws.Range("I21").Value = 1.192514257685282
$ws.Range("J21").Value = 0.1072694729077046
$ws.Range("K21").Value = 0.8426316391769149
$ws.Range("M21").Value = 0.5018526938483916
$ws.Range("N21").Value = 2.425617509389781
# Row 22
$ws.Range("B22").Value = 0.7066485769352937
$ws.Range("C22").Value = 0.1919160789837804
$ws.Range("E22").Value = 0.2640468577629846
$ws.Range("F22").Value = 3.064272440401425
$ws.Range("G22").Value = 0.002505576602344369
$ws.Range("I22").Value = 1.19450151810161
$ws.Range("J22").Value = 0.1067810324432816
$ws.Range("K22").Value = 0.8817475609874066
$ws.Range("M22").Value = 0.5182801958927143
$ws.Range("N22").Value = 2.413587362196424
# Row 23
$ws.Range("B23").Value = 0.6894037662612504
$ws.Range("C23").Value = 0.1877405051631342
$ws.Range("E23").Value = 0.2613488717560628
$ws.Range("F23").Value = 3.048215854169229
$ws.Range("G23").Value = 0.002506949048063447
$ws.Range("I23").Value = 1.193397778937829
$ws.Range("J23").Value = 0.107039124793916
$ws.Range("K23").Value = 0.8608410332613232
$ws.Range("M23").Value = 0.5094936845521261
$ws.Range("N23").Value = 2.419949044757956
# Row 24
$ws.Range("B24").Value = 0.6245781080722566
$ws.Range("C24").Value = 0.1720938378878714
$ws.Range("E24").Value = 0.2513300931742677
$ws.Range("F24").Value = 2.99001155403127
$ws.Range("G24").Value = 0.00251234703038974
$ws.Range("I24").Value = 1.19016335364271
$ws.Range("J24").Value = 0.1080736328684253
$ws.Range("K24").Value = 0.7823110185043163
$ws.Range("M24").Value = 0.4766302384649919
$ws.Range("N24").Value = 2.445338354964456
# Row 25
$ws.Range("B25").Value = 0.5556930632583317
$ws.Range("C25").Value = 0.1555703080309172
$ws.Range("E25").Value = 0.2409506533392047
$ws.Range("F25").Value = 2.932786299594667
$ws.Range("G25").Value = 0.002518600109588598
$ws.Range("I25").Value = 1.188688399487575
$ws.Range("J25").Value = 0.1093110928707617
$ws.Range("K25").Value = 0.6989884673415361
$ws.Range("M25").Value = 0.4420668625568638
$ws.Range("N25").Value = 2.475483316348274

Write-Host "Updated 380 kV case values"